# "Generate Report for Archive"
# The localization status report is regenerated: every cell that held the
# "Ready for handoff" status moves on to "In Translation", and the Status
# column on each sheet is narrowed to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn / de-de status columns (E, F) ---
foreach ($r in 2..4) {
    if ($overview.Cells.Item($r, 5).Value2 -eq $oldStatus) {
        $overview.Cells.Item($r, 5).Value = $newStatus
    }
    if ($overview.Cells.Item($r, 6).Value2 -eq $oldStatus) {
        $overview.Cells.Item($r, 6).Value = $newStatus
    }
}

# --- zh-cn / de-de sheets: Status column (C) ---
foreach ($r in 2..4) {
    if ($zhcn.Cells.Item($r, 3).Value2 -eq $oldStatus) {
        $zhcn.Cells.Item($r, 3).Value = $newStatus
    }
    if ($dede.Cells.Item($r, 3).Value2 -eq $oldStatus) {
        $dede.Cells.Item($r, 3).Value = $newStatus
    }
}

# --- Re-fit the Status column widths now that the text is shorter ---
# (Target authored width is 13.4101845877511 characters; the runtime quantizes
# ColumnWidth writes to the nearest 1/6 character, so 12.5 is the input that
# lands closest to that authored width.)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
